# Scheduled-runner price refresh for the Diabolos Profits workbook.
# Updates cached market-board figures (currentAveragePrice* / LevePrice* /
# LeveProfit*) in columns H:N on each crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Cells.Item(98, 8).Value = 3760.4348
$ws.Cells.Item(98, 9).Value = 2795.2354
$ws.Cells.Item(98, 10).Value = 6495.1665
$ws.Cells.Item(98, 11).Value = 2795.2354
$ws.Cells.Item(98, 12).Value = 6495.1665
$ws.Cells.Item(98, 13).Value = -1297.2354
$ws.Cells.Item(98, 14).Value = -9491.166499999999

# Row 113: Amaro Kart / Starch Glue
$ws.Cells.Item(113, 8).Value = 166670670
$ws.Cells.Item(113, 9).Value = 333337000
$ws.Cells.Item(113, 10).Value = 4333.3335
$ws.Cells.Item(113, 11).Value = 333337000
$ws.Cells.Item(113, 12).Value = 4333.3335
$ws.Cells.Item(113, 13).Value = -333333746
$ws.Cells.Item(113, 14).Value = -10841.3335

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Cells.Item(122, 8).Value = 3760.4348
$ws.Cells.Item(122, 9).Value = 2795.2354
$ws.Cells.Item(122, 10).Value = 6495.1665
$ws.Cells.Item(122, 11).Value = 8385.706200000001
$ws.Cells.Item(122, 12).Value = 19485.4995
$ws.Cells.Item(122, 13).Value = -5935.706200000001
$ws.Cells.Item(122, 14).Value = -24385.4995

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 4298.4062
$ws.Cells.Item(132, 9).Value = 4185.6665
$ws.Cells.Item(132, 10).Value = 4636.625
$ws.Cells.Item(132, 11).Value = 12556.9995
$ws.Cells.Item(132, 12).Value = 13909.875
$ws.Cells.Item(132, 13).Value = -10026.9995
$ws.Cells.Item(132, 14).Value = -18969.875

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 4791.64
$ws.Cells.Item(138, 9).Value = 5339
$ws.Cells.Item(138, 10).Value = 4618.7896
$ws.Cells.Item(138, 11).Value = 16017
$ws.Cells.Item(138, 12).Value = 13856.3688
$ws.Cells.Item(138, 13).Value = -10877
$ws.Cells.Item(138, 14).Value = -24136.3688

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 8785.237999999999
$ws.Cells.Item(141, 9).Value = 10766.667
$ws.Cells.Item(141, 10).Value = 3831.6667
$ws.Cells.Item(141, 11).Value = 32300.001
$ws.Cells.Item(141, 12).Value = 11495.0001
$ws.Cells.Item(141, 13).Value = -27120.001

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2618.2104
$ws.Cells.Item(61, 9).Value = 2298.394
$ws.Cells.Item(61, 10).Value = 4729
$ws.Cells.Item(61, 11).Value = 2298.394
$ws.Cells.Item(61, 12).Value = 4729
$ws.Cells.Item(61, 13).Value = -2086.394

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 4182.212
$ws.Cells.Item(74, 9).Value = 3183.64
$ws.Cells.Item(74, 10).Value = 7302.75
$ws.Cells.Item(74, 11).Value = 3183.64
$ws.Cells.Item(74, 12).Value = 7302.75
$ws.Cells.Item(74, 13).Value = -2309.64
$ws.Cells.Item(74, 14).Value = -9050.75

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 4182.212
$ws.Cells.Item(77, 9).Value = 3183.64
$ws.Cells.Item(77, 10).Value = 7302.75
$ws.Cells.Item(77, 11).Value = 15918.2
$ws.Cells.Item(77, 12).Value = 36513.75
$ws.Cells.Item(77, 13).Value = -11550.2
$ws.Cells.Item(77, 14).Value = -45249.75

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 3080633.8
$ws.Cells.Item(102, 9).Value = 3850167.2
$ws.Cells.Item(102, 10).Value = 2500
$ws.Cells.Item(102, 11).Value = 3850167.2
$ws.Cells.Item(102, 12).Value = 2500
$ws.Cells.Item(102, 13).Value = -3848545.2

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 3950.1724
$ws.Cells.Item(132, 9).Value = 3482.8
$ws.Cells.Item(132, 10).Value = 6871.25
$ws.Cells.Item(132, 11).Value = 10448.4
$ws.Cells.Item(132, 12).Value = 20613.75
$ws.Cells.Item(132, 13).Value = -7918.400000000001

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2618.2104
$ws.Cells.Item(136, 9).Value = 2298.394
$ws.Cells.Item(136, 10).Value = 4729
$ws.Cells.Item(136, 11).Value = 6895.181999999999
$ws.Cells.Item(136, 12).Value = 14187
$ws.Cells.Item(136, 13).Value = -4345.181999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal / High Steel Nugget
$ws.Cells.Item(94, 8).Value = 2268.1333
$ws.Cells.Item(94, 9).Value = 2386.3076
$ws.Cells.Item(94, 10).Value = 1500
$ws.Cells.Item(94, 11).Value = 2386.3076
$ws.Cells.Item(94, 12).Value = 1500
$ws.Cells.Item(94, 13).Value = -1935.3076
$ws.Cells.Item(94, 14).Value = -2402

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 5401.357
$ws.Cells.Item(99, 9).Value = 4556.4546
$ws.Cells.Item(99, 10).Value = 8499.333000000001
$ws.Cells.Item(99, 11).Value = 4556.4546
$ws.Cells.Item(99, 12).Value = 8499.333000000001
$ws.Cells.Item(99, 13).Value = -3058.4546
$ws.Cells.Item(99, 14).Value = -11495.333

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 29531658
$ws.Cells.Item(107, 9).Value = 168754.42
$ws.Cells.Item(107, 10).Value = 100002630
$ws.Cells.Item(107, 11).Value = 168754.42
$ws.Cells.Item(107, 12).Value = 100002630
$ws.Cells.Item(107, 13).Value = -166834.42

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 3680.3777
$ws.Cells.Item(31, 9).Value = 2200.4707
$ws.Cells.Item(31, 10).Value = 4578.893
$ws.Cells.Item(31, 11).Value = 2200.4707
$ws.Cells.Item(31, 12).Value = 4578.893
$ws.Cells.Item(31, 13).Value = -1905.4707

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 3680.3777
$ws.Cells.Item(34, 9).Value = 2200.4707
$ws.Cells.Item(34, 10).Value = 4578.893
$ws.Cells.Item(34, 11).Value = 2200.4707
$ws.Cells.Item(34, 12).Value = 4578.893
$ws.Cells.Item(34, 13).Value = -1998.4707

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 438079.3
$ws.Cells.Item(132, 9).Value = 3744.4167
$ws.Cells.Item(132, 10).Value = 911899.2
$ws.Cells.Item(132, 11).Value = 11233.2501
$ws.Cells.Item(132, 12).Value = 2735697.6
$ws.Cells.Item(132, 13).Value = -8703.250100000001
$ws.Cells.Item(132, 14).Value = -2740757.6

$ws = $wb.Worksheets.Item("CUL")
# Row 81: It Goes Down Smoothly / Frozen Spirits
$ws.Cells.Item(81, 8).Value = 11877.75
$ws.Cells.Item(81, 9).Value = 1000
$ws.Cells.Item(81, 10).Value = 13431.714
$ws.Cells.Item(81, 11).Value = 3000
$ws.Cells.Item(81, 12).Value = 40295.142
$ws.Cells.Item(81, 13).Value = -1877
$ws.Cells.Item(81, 14).Value = -42541.142

# Row 84: Quenching the Flame (L) / Frozen Spirits
$ws.Cells.Item(84, 8).Value = 11877.75
$ws.Cells.Item(84, 9).Value = 1000
$ws.Cells.Item(84, 10).Value = 13431.714
$ws.Cells.Item(84, 11).Value = 9000
$ws.Cells.Item(84, 12).Value = 120885.426
$ws.Cells.Item(84, 13).Value = -3384
$ws.Cells.Item(84, 14).Value = -132117.426

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 12639.19
$ws.Cells.Item(131, 9).Value = 2025.5
$ws.Cells.Item(131, 10).Value = 16884.666
$ws.Cells.Item(131, 11).Value = 6076.5
$ws.Cells.Item(131, 12).Value = 50653.99800000001
$ws.Cells.Item(131, 13).Value = -1036.5

# Row 132: More Mezcal / Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 2932.8547
$ws.Cells.Item(132, 9).Value = 1279.6666
$ws.Cells.Item(132, 10).Value = 3016.9153
$ws.Cells.Item(132, 11).Value = 11516.9994
$ws.Cells.Item(132, 12).Value = 27152.2377
$ws.Cells.Item(132, 13).Value = -8986.999400000001
$ws.Cells.Item(132, 14).Value = -32212.2377

$ws = $wb.Worksheets.Item("GSM")
# Row 51: When We Were Blings / Mythril Ear Cuffs
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Cells.Item(70, 8).Value = 14284
$ws.Cells.Item(70, 9).Value = 11997.75
$ws.Cells.Item(70, 10).Value = 17332.334
$ws.Cells.Item(70, 11).Value = 11997.75
$ws.Cells.Item(70, 12).Value = 17332.334
$ws.Cells.Item(70, 13).Value = -11727.75
$ws.Cells.Item(70, 14).Value = -17872.334

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Cells.Item(73, 8).Value = 14284
$ws.Cells.Item(73, 9).Value = 11997.75
$ws.Cells.Item(73, 10).Value = 17332.334
$ws.Cells.Item(73, 11).Value = 11997.75
$ws.Cells.Item(73, 12).Value = 17332.334
$ws.Cells.Item(73, 13).Value = -11061.75
$ws.Cells.Item(73, 14).Value = -19204.334

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Cells.Item(97, 8).Value = 755
$ws.Cells.Item(97, 9).Value = 818.0476
$ws.Cells.Item(97, 10).Value = 644.6667
$ws.Cells.Item(97, 11).Value = 818.0476
$ws.Cells.Item(97, 12).Value = 644.6667
$ws.Cells.Item(97, 13).Value = -322.0476
$ws.Cells.Item(97, 14).Value = -1636.6667

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Cells.Item(122, 8).Value = 507741.38
$ws.Cells.Item(122, 9).Value = 696706.9
$ws.Cells.Item(122, 10).Value = 3833.3333
$ws.Cells.Item(122, 11).Value = 2090120.7
$ws.Cells.Item(122, 12).Value = 11499.9999
$ws.Cells.Item(122, 13).Value = -2087670.7

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value = 4074.88
$ws.Cells.Item(40, 9).Value = 3630.6155
$ws.Cells.Item(40, 10).Value = 4556.1665
$ws.Cells.Item(40, 11).Value = 3630.6155
$ws.Cells.Item(40, 12).Value = 4556.1665
$ws.Cells.Item(40, 13).Value = -3494.6155
$ws.Cells.Item(40, 14).Value = -4828.1665

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Cells.Item(55, 8).Value = 1473.375
$ws.Cells.Item(55, 9).Value = 698.5
$ws.Cells.Item(55, 10).Value = 2248.25
$ws.Cells.Item(55, 11).Value = 698.5
$ws.Cells.Item(55, 12).Value = 2248.25
$ws.Cells.Item(55, 13).Value = -525.5
$ws.Cells.Item(55, 14).Value = -2594.25

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Cells.Item(100, 8).Value = 4687
$ws.Cells.Item(100, 9).Value = 4449.5
$ws.Cells.Item(100, 10).Value = 4924.5
$ws.Cells.Item(100, 11).Value = 4449.5
$ws.Cells.Item(100, 12).Value = 4924.5
$ws.Cells.Item(100, 13).Value = -3908.5
$ws.Cells.Item(100, 14).Value = -6006.5

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Cells.Item(136, 8).Value = 7007.478
$ws.Cells.Item(136, 9).Value = 7119.7896
$ws.Cells.Item(136, 10).Value = 6474
$ws.Cells.Item(136, 11).Value = 21359.3688
$ws.Cells.Item(136, 12).Value = 19422
$ws.Cells.Item(136, 13).Value = -18809.3688

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Cells.Item(62, 8).Value = 3657076.8
$ws.Cells.Item(62, 9).Value = 8765084
$ws.Cells.Item(62, 10).Value = 8500
$ws.Cells.Item(62, 11).Value = 8765084
$ws.Cells.Item(62, 12).Value = 8500
$ws.Cells.Item(62, 13).Value = -8764460

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Cells.Item(65, 8).Value = 3657076.8
$ws.Cells.Item(65, 9).Value = 8765084
$ws.Cells.Item(65, 10).Value = 8500
$ws.Cells.Item(65, 11).Value = 43825420
$ws.Cells.Item(65, 12).Value = 42500
$ws.Cells.Item(65, 13).Value = -43822300

# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Cells.Item(81, 8).Value = 14291670
$ws.Cells.Item(81, 9).Value = 4849.4443
$ws.Cells.Item(81, 10).Value = 40007944
$ws.Cells.Item(81, 11).Value = 9698.8886
$ws.Cells.Item(81, 12).Value = 80015888
$ws.Cells.Item(81, 13).Value = -8637.8886

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Cells.Item(84, 8).Value = 14291670
$ws.Cells.Item(84, 9).Value = 4849.4443
$ws.Cells.Item(84, 10).Value = 40007944
$ws.Cells.Item(84, 11).Value = 48494.443
$ws.Cells.Item(84, 12).Value = 400079440
$ws.Cells.Item(84, 13).Value = -43190.443

# Row 113: A Tender Table / Pixie Floss
$ws.Cells.Item(113, 8).Value = 6344.88
$ws.Cells.Item(113, 9).Value = 7436.875
$ws.Cells.Item(113, 10).Value = 4403.5557
$ws.Cells.Item(113, 11).Value = 22310.625
$ws.Cells.Item(113, 12).Value = 13210.6671
$ws.Cells.Item(113, 13).Value = -20140.625

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Cells.Item(122, 8).Value = 1944.9722
$ws.Cells.Item(122, 9).Value = 1371.7
$ws.Cells.Item(122, 10).Value = 2661.5625
$ws.Cells.Item(122, 11).Value = 4115.1
$ws.Cells.Item(122, 12).Value = 7984.6875
$ws.Cells.Item(122, 13).Value = -1665.1

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 640912.9
$ws.Cells.Item(132, 9).Value = 1184973.1
$ws.Cells.Item(132, 10).Value = 24311.334
$ws.Cells.Item(132, 11).Value = 3554919.3
$ws.Cells.Item(132, 12).Value = 72934.00199999999
$ws.Cells.Item(132, 13).Value = -3552389.3
